# Users_Info.xlsx update:
#  - Venues sheet: rename header to "Venues_ID", trim/replace the list of venue
#    IDs (drop 117/118/119/217/218/219, add 311-315), shrink the table/range
#    from A7:A25 down to A1:A18, restyle the table.
#  - Add a new "Issues" sheet with an Issue_ID/Description table.
#  - Make "Issues" the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Venues sheet: update header + data, resize table, restyle
# ---------------------------------------------------------------------------
$venues = $wb.Worksheets.Item("Venues")

# Resize the table first (while the header cell still reads "Venues") so the
# ListObject's column keeps following the header cell instead of renaming it.
$venuesTable = $venues.ListObjects.Item(1)
$venuesTable.Resize($venues.Range("A1:A18"))

$venues.Range("A1").Value = "Venues_ID"

$venueIds = @(111,112,113,114,115,116,211,212,213,214,215,216,311,312,313,314,315)
$venueArr = New-Object 'object[,]' $venueIds.Count,1
for ($i = 0; $i -lt $venueIds.Count; $i++) {
    $venueArr[$i,0] = $venueIds[$i]
}
$venues.Range("A2:A18").Value = $venueArr

# Drop the now-unused trailing rows (old sheet ran to row 25).
$venues.Range("A19:A25").ClearContents()

$venuesTable.TableStyle = "TableStyleMedium7"

# ---------------------------------------------------------------------------
# 2. Add the new "Issues" sheet after "Venues"
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$issues = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$issues.Name = "Issues"

# Match the original authoring order: IDs first, then the two headers, then
# the descriptions (keeps the shared-string table in the same append order).
$issueIds = @("I001","I002","I003","I004","I005","I006","I007")
$issueDescriptions = @(
    "Failure of Login Account",
    "Failure of Connect WiFi",
    "Mic Sound Weak",
    "Failure of Operate Smart TV",
    "Sound Effect of Playing Video Weak",
    "Failure of Project Videos",
    "Others"
)

for ($i = 0; $i -lt $issueIds.Count; $i++) {
    $issues.Range("A$($i + 2)").Value = $issueIds[$i]
}

$issues.Range("A1").Value = "Issue_ID"
$issues.Range("B1").Value = "Description"

for ($i = 0; $i -lt $issueDescriptions.Count; $i++) {
    $issues.Range("B$($i + 2)").Value = $issueDescriptions[$i]
}

$issuesTable = $issues.ListObjects.Add(1, $issues.Range("A1:B8"), [System.Reflection.Missing]::Value, 1)
$issuesTable.Name = "表格5"
$issuesTable.TableStyle = "TableStyleMedium7"

$issues.Columns.Item(2).ColumnWidth = 40.18

# ---------------------------------------------------------------------------
# 3. Selections / active sheet to match the edited file
# ---------------------------------------------------------------------------
$venues.Activate()
$venues.Range("E4").Select()

$issues.Activate()
$issues.Range("D3").Select()
